# Delete slides and then update shapes on the remaining slide.
$p = $ppt.ActivePresentation

# Delete slides 4 and 3 (in descending order so indices stay valid).
$p.Slides.Item(4).Delete()
$p.Slides.Item(3).Delete()

# Remaining slide 2 ("Shape 129" deck) needs shape repositioning/resizing.
$s = $p.Slides.Item(2)

# Shape id=133 "Google Shape;133;ge6d407e7a1_0_85" - textbox with the main body text.
$body = $s.Shapes.Item(3)
$body.Left = 508375
$body.Top = 1032564
$body.Width = 8386475
$body.Height = 5740003

# Add three extra blank paragraphs after the 3rd paragraph (the one
# containing the standalone <a:endParaRPr sz="1700" dirty="0"/>).
$tr = $body.TextFrame.TextRange
$tr.Paragraphs(3).InsertAfter("`r`v`r`v`r")

# Shape id=2 "Google Shape;120;ge6d407e7a1_0_27" - black chmod box.
$box = $s.Shapes.Item(6)
$box.Top = 4101403
